# Daily update: insert a new "today" row (2021/12/10) at the top of the
# history table on every sheet, pushing all existing rows down by one.
#
# All five sheets in this workbook are rolling daily logs with newest date
# in row 2 (row 1 is the header). The commit simply prepends one more day
# of data to each log.

$wb = $excel.ActiveWorkbook

function Set-TextCell($cell, [string]$text) {
    # Force literal-text storage so values that look numeric/date-like
    # ("202201", "2021/12/10") are not auto-coerced by Excel into a
    # number or a date serial.
    $cell.NumberFormat = "@"
    $cell.Value = $text
}

# ---------------------------------------------------------------------
# Sheet 1: 台指期換倉成本計算  (A:日期 B:月份 C:結算價 D:未沖銷契約量 E:金額 F:成本)
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)
$ws1.Rows.Item(2).Insert()
Set-TextCell $ws1.Cells.Item(2,1) "日期：2021/12/10"
Set-TextCell $ws1.Cells.Item(2,2) "202201"
$ws1.Cells.Item(2,3).Value = 17763
$ws1.Cells.Item(2,4).Value = 13989
$ws1.Cells.Item(2,5).Value = 32008926
$ws1.Cells.Item(2,6).Value = 17676

# ---------------------------------------------------------------------
# Sheet 2: 散戶多空力道  (A:日期 B:散戶多空力道)
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item(2)
$ws2.Rows.Item(2).Insert()
Set-TextCell $ws2.Cells.Item(2,1) "日期：2021/12/10"
$ws2.Cells.Item(2,2).Value = 0.06

# ---------------------------------------------------------------------
# Sheet 3: 三大法人買賣金額  (A:日期 B:外資 C:內資)
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item(3)
$ws3.Rows.Item(2).Insert()
Set-TextCell $ws3.Cells.Item(2,1) "110年12月10日"
$ws3.Cells.Item(2,2).Value = -32.7
$ws3.Cells.Item(2,3).Value = -55.16

# ---------------------------------------------------------------------
# Sheet 4: 大盤多空點位  (A:日期 B:隔日多空點位)
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)
$ws4.Rows.Item(2).Insert()
Set-TextCell $ws4.Cells.Item(2,1) "110年12月10日"
$ws4.Cells.Item(2,2).Value = 17832.17

# ---------------------------------------------------------------------
# Sheet 5: 期貨大額交易人未沖銷部位  (A:日期 B..N: positions)
# ---------------------------------------------------------------------
$ws5 = $wb.Worksheets.Item(5)
$ws5.Rows.Item(2).Insert()
Set-TextCell $ws5.Cells.Item(2,1) "2021/12/10"
$ws5.Cells.Item(2,2).Value  = 48395
$ws5.Cells.Item(2,3).Value  = 57141
$ws5.Cells.Item(2,4).Value  = 944
$ws5.Cells.Item(2,5).Value  = 879
$ws5.Cells.Item(2,6).Value  = 24508
$ws5.Cells.Item(2,7).Value  = 50924
$ws5.Cells.Item(2,8).Value  = -1296
$ws5.Cells.Item(2,9).Value  = 520
$ws5.Cells.Item(2,10).Value = -26416
$ws5.Cells.Item(2,11).Value = -1816
$ws5.Cells.Item(2,12).Value = 2240
$ws5.Cells.Item(2,13).Value = 359
$ws5.Cells.Item(2,14).Value = 1881
